$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.087693701173381
$ws.Range("D2").Value = 1.090087753345538
$ws.Range("E2").Value = 1.090023313609521
$ws.Range("F2").Value = 1.10124700516267
$ws.Range("I2").Value = 1.069491623842711
$ws.Range("J2").Value = 1.092537525210106
$ws.Range("K2").Value = 1.092736205283928
$ws.Range("L2").Value = 1.092671930797378
$ws.Range("M2").Value = 1.10386715452933

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.089082650834224
$ws.Range("D3").Value = 1.091243477302079
$ws.Range("E3").Value = 1.091276826003912
$ws.Range("F3").Value = 1.102559642743688
$ws.Range("I3").Value = 1.070024886520649
$ws.Range("J3").Value = 1.093588501226809
$ws.Range("K3").Value = 1.093711178777864
$ws.Range("L3").Value = 1.093744448162252
$ws.Range("M3").Value = 1.105000724982897

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.089980194183307
$ws.Range("D4").Value = 1.091990209655082
$ws.Range("E4").Value = 1.092086550247943
$ws.Range("F4").Value = 1.10340800427791
$ws.Range("I4").Value = 1.070368042182545
$ws.Range("J4").Value = 1.09426682930401
$ws.Range("K4").Value = 1.09434036633362
$ws.Range("L4").Value = 1.094436489356092
$ws.Range("M4").Value = 1.105732658577304

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.09035723966544
$ws.Range("D5").Value = 1.092303876991763
$ws.Range("E5").Value = 1.092426632239524
$ws.Range("F5").Value = 1.103764420108821
$ws.Range("I5").Value = 1.070511852139323
$ws.Range("J5").Value = 1.094551589998476
$ws.Range("K5").Value = 1.094604477091012
$ws.Range("L5").Value = 1.094726961304534
$ws.Range("M5").Value = 1.106039993833351

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.090420530899811
$ws.Range("D6").Value = 1.092356528048975
$ws.Range("E6").Value = 1.092483714519825
$ws.Range("F6").Value = 1.103824250229934
$ws.Range("I6").Value = 1.070535971988042
$ws.Range("J6").Value = 1.094599378744004
$ws.Range("K6").Value = 1.094648799121887
$ws.Range("L6").Value = 1.094775705851035
$ws.Range("M6").Value = 1.106091575227242

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.089985233383665
$ws.Range("D7").Value = 1.091994401905894
$ws.Range("E7").Value = 1.092091095715464
$ws.Range("F7").Value = 1.103412767642631
$ws.Range("I7").Value = 1.070369965553977
$ws.Range("J7").Value = 1.094270635889946
$ws.Range("K7").Value = 1.094343896960257
$ws.Range("L7").Value = 1.094440372466913
$ws.Range("M7").Value = 1.105736766652964

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.088163355009003
$ws.Range("D8").Value = 1.09047856549578
$ws.Range("E8").Value = 1.090447233056691
$ws.Range("F8").Value = 1.101690827643009
$ws.Range("I8").Value = 1.06967223751798
$ws.Range("J8").Value = 1.092893066923423
$ws.Range("K8").Value = 1.093066053436525
$ws.Range("L8").Value = 1.093034799375434
$ws.Range("M8").Value = 1.104250576100673

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.084943533962803
$ws.Range("D9").Value = 1.087798876399322
$ws.Range("E9").Value = 1.087539736383975
$ws.Range("F9").Value = 1.098648658097259
$ws.Range("I9").Value = 1.068428084321367
$ws.Range("J9").Value = 1.090452223236229
$ws.Range("K9").Value = 1.090801249372902
$ws.Range("L9").Value = 1.09054286947142
$ws.Range("M9").Value = 1.101619556008169

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.082790291539147
$ws.Range("D10").Value = 1.086006386852132
$ws.Range("E10").Value = 1.085593840315151
$ws.Range("F10").Value = 1.096614941751933
$ws.Range("I10").Value = 1.067588634806157
$ws.Range("J10").Value = 1.088815723389464
$ws.Range("K10").Value = 1.089282343476879
$ws.Range("L10").Value = 1.088871126390436
$ws.Range("M10").Value = 1.099857087590932

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.081856247571994
$ws.Range("D11").Value = 1.085228731714013
$ws.Range("E11").Value = 1.084749383451912
$ws.Range("F11").Value = 1.095732924199066
$ws.Range("I11").Value = 1.067222734341838
$ws.Range("J11").Value = 1.088104841954761
$ws.Range("K11").Value = 1.088622442699819
$ws.Range("L11").Value = 1.088144699399798
$ws.Range("M11").Value = 1.099091851500841

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.081509043461104
$ws.Range("D12").Value = 1.084939646304385
$ws.Range("E12").Value = 1.084435427519776
$ws.Range("F12").Value = 1.095405086678062
$ws.Range("I12").Value = 1.067086457022162
$ws.Range("J12").Value = 1.087840443207345
$ws.Range("K12").Value = 1.088376990194113
$ws.Range("L12").Value = 1.087874483210817
$ws.Range("M12").Value = 1.09880729144066

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.081583531733422
$ws.Range("D13").Value = 1.08500166657975
$ws.Range("E13").Value = 1.084502785218742
$ws.Range("F13").Value = 1.09547541885009
$ws.Range("I13").Value = 1.067115705562607
$ws.Range("J13").Value = 1.08789717336603
$ws.Range("K13").Value = 1.088429655870356
$ws.Range("L13").Value = 1.08793246318719
$ws.Range("M13").Value = 1.098868344986311

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.081827552886724
$ws.Range("D14").Value = 1.085204840538089
$ws.Range("E14").Value = 1.084723437646997
$ws.Range("F14").Value = 1.095705829506597
$ws.Range("I14").Value = 1.067211477091695
$ws.Range("J14").Value = 1.088082993755158
$ws.Range("K14").Value = 1.088602160400639
$ws.Range("L14").Value = 1.0881223712103
$ws.Range("M14").Value = 1.099068336182132

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.081977867902008
$ws.Range("D15").Value = 1.085329992135978
$ws.Range("E15").Value = 1.084859350674374
$ws.Range("F15").Value = 1.095847764195562
$ws.Range("I15").Value = 1.067270436555803
$ws.Range("J15").Value = 1.088197437859191
$ws.Range("K15").Value = 1.08870840147183
$ws.Range("L15").Value = 1.088239328116086
$ws.Range("M15").Value = 1.099191515148508

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.08285224507252
$ws.Range("D16").Value = 1.08605796533074
$ws.Range("E16").Value = 1.085649844208364
$ws.Range("F16").Value = 1.096673448252685
$ws.Range("I16").Value = 1.067612867311858
$ws.Range("J16").Value = 1.088862854054631
$ws.Range("K16").Value = 1.089326092065072
$ws.Range("L16").Value = 1.088919282731688
$ws.Range("M16").Value = 1.099907829593195

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.083400265309028
$ws.Range("D17").Value = 1.086514199813275
$ws.Range("E17").Value = 1.086145194799398
$ws.Range("F17").Value = 1.097190997781441
$ws.Range("I17").Value = 1.067827016842273
$ws.Range("J17").Value = 1.089279641624416
$ws.Range("K17").Value = 1.089712959486814
$ws.Range("L17").Value = 1.08934511334513
$ws.Range("M17").Value = 1.100356595166449

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.083719755015803
$ws.Range("D18").Value = 1.086780169730763
$ws.Range("E18").Value = 1.086433944296535
$ws.Range("F18").Value = 1.097492740505547
$ws.Range("I18").Value = 1.067951693999116
$ws.Range("J18").Value = 1.089522528464413
$ws.Range("K18").Value = 1.089938400330993
$ws.Range("L18").Value = 1.089593247361876
$ws.Range("M18").Value = 1.10061815263569

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.083828665608352
$ws.Range("D19").Value = 1.086870834335771
$ws.Range("E19").Value = 1.086532370040521
$ws.Range("F19").Value = 1.097595604227403
$ws.Range("I19").Value = 1.06799416631878
$ws.Range("J19").Value = 1.089605309731287
$ws.Range("K19").Value = 1.090015233953674
$ws.Range("L19").Value = 1.089677813129164
$ws.Range("M19").Value = 1.100707303231884

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.083341484651404
$ws.Range("D20").Value = 1.086465265096751
$ws.Range("E20").Value = 1.086092067044367
$ws.Range("F20").Value = 1.097135483639102
$ws.Range("I20").Value = 1.067804064710208
$ws.Range("J20").Value = 1.089234946881589
$ws.Range("K20").Value = 1.089671474276733
$ws.Range("L20").Value = 1.089299451213747
$ws.Range("M20").Value = 1.100308467584739

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.081755701932602
$ws.Range("D21").Value = 1.085145017251646
$ws.Range("E21").Value = 1.084658468978577
$ws.Range("F21").Value = 1.095637985346603
$ws.Range("I21").Value = 1.067183284880049
$ws.Range("J21").Value = 1.088028283868271
$ws.Range("K21").Value = 1.088551371429081
$ws.Range("L21").Value = 1.088066458810188
$ws.Range("M21").Value = 1.099009452562811

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.080757157847673
$ws.Range("D22").Value = 1.084313593747887
$ws.Range("E22").Value = 1.083755445196545
$ws.Range("F22").Value = 1.09469519071955
$ws.Range("I22").Value = 1.066790859226767
$ws.Range("J22").Value = 1.087267604020678
$ws.Range("K22").Value = 1.087845172077838
$ws.Range("L22").Value = 1.087288974939132
$ws.Range("M22").Value = 1.098190871706243

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.081286649388855
$ws.Range("D23").Value = 1.084754474949073
$ws.Range("E23").Value = 1.08423431489266
$ws.Range("F23").Value = 1.095195105197545
$ws.Range("I23").Value = 1.066999093147652
$ws.Range("J23").Value = 1.087671046417909
$ws.Range("K23").Value = 1.088219727797429
$ws.Range("L23").Value = 1.087701349238956
$ws.Range("M23").Value = 1.098624993012474

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.083368045598098
$ws.Range("D24").Value = 1.086487377033221
$ws.Range("E24").Value = 1.086116073746175
$ws.Range("F24").Value = 1.097160568507706
$ws.Range("I24").Value = 1.067814436509346
$ws.Range("J24").Value = 1.089255143186168
$ws.Range("K24").Value = 1.089690220313363
$ws.Range("L24").Value = 1.089320084724451
$ws.Range("M24").Value = 1.100330214984584

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.08577708922158
$ws.Range("D25").Value = 1.088492684425884
$ws.Range("E25").Value = 1.088292705006486
$ws.Range("F25").Value = 1.099436097536459
$ws.Range("I25").Value = 1.068751481545062
$ws.Range("J25").Value = 1.091084854174532
$ws.Range("K25").Value = 1.091388330432065
$ws.Range("L25").Value = 1.091188915224357
$ws.Range("M25").Value = 1.102301207566203
